$wb = $excel.ActiveWorkbook

# --- Sheet1 ---
$ws1 = $wb.Worksheets.Item("Sheet1")

# Add "Barrier" to the Defence/Type list (column I) at row 9
$ws1.Range("I9").Value = "Barrier"

# Remove "Barrier" from the Spell Type list (column D) at row 11
$ws1.Range("D11").ClearContents()

# Update the active selection to I10 as in the edited file
$ws1.Range("I10").Select()

# --- Sheet2 ---
$ws2 = $wb.Worksheets.Item("Sheet2")

# Attack section: Resource -> Hand Seal, Spell Type -> Barrier
$ws2.Range("B3").Value = "Hand Seal"
$ws2.Range("C3").Value = "Barrier"

# Defence section: Type -> Barrier, Number -> Roll Count
$ws2.Range("B5").Value = "Barrier"
$ws2.Range("C5").Value = "Roll Count"

# Update the active selection to E9 as in the edited file
$ws2.Range("E9").Select()
$ws2.Activate()
